$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.583.28"
$ws.Range("E2").Value = "  +0.36%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.692.62"
$ws.Range("E3").Value = "  +0.27%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "676.29"
$ws.Range("E5").Value = "  -1.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "161.40"
$ws.Range("E6").Value = "  +1.38%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +0.55%  "

$ws.Range("E9").Value = "  +1.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.12"
$ws.Range("E10").Value = "  -0.37%  "

$ws.Range("E11").Value = "  +1.46%  "

$ws.Range("E12").Value = "  +0.67%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "32.51"
$ws.Range("E13").Value = "  +0.64%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.686.62"
$ws.Range("E14").Value = "  +0.04%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "69.507.46"
$ws.Range("E15").Value = "  +0.28%  "

$ws.Range("E16").Value = "  +2.31%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "16.02"
$ws.Range("E17").Value = "  +1.09%  "

$ws.Range("E18").Value = "  +0.56%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "471.37"
$ws.Range("E19").Value = "  +1.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.83"
$ws.Range("E20").Value = "  -2.40%  "

$ws.Range("E21").Value = "  +0.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "80.56"
$ws.Range("E22").Value = "  +1.60%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.838.84"
$ws.Range("E23").Value = "  +0.23%  "

$ws.Range("E24").Value = "  -0.11%  "

$ws.Range("E25").Value = "  +1.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.88"
$ws.Range("E26").Value = "  -0.41%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.12"
$ws.Range("E27").Value = "  -0.39%  "

$ws.Range("E28").Value = "  -0.19%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.75"
$ws.Range("E29").Value = "  +0.77%  "

$ws.Range("E30").Value = "  +0.57%  "

$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("E32").Value = "  -0.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.98"
$ws.Range("E33").Value = "  +1.10%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.682.31"
$ws.Range("E34").Value = "  +0.68%  "

$ws.Range("E35").Value = "  +1.06%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.46"
$ws.Range("E36").Value = "  +4.04%  "

$ws.Range("E37").Value = "  +1.66%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.25"
$ws.Range("E39").Value = "  -0.90%  "

$ws.Range("E40").Value = "  -0.03%  "

$ws.Range("E41").Value = "  +0.44%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "167.30"
$ws.Range("E42").Value = "  +0.32%  "

$ws.Range("E43").Value = "  +0.13%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "46.50"
$ws.Range("E44").Value = "  -2.48%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.76"
$ws.Range("E45").Value = "  +1.17%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000280"
$ws.Range("E46").Value = "  +1.93%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "27.93"
$ws.Range("E47").Value = "  -1.11%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.30"
$ws.Range("E48").Value = "  -0.87%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.10"
$ws.Range("E49").Value = "  -1.38%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.88"
$ws.Range("E50").Value = "  +1.47%  "

$ws.Range("E51").Value = "  +1.91%  "
